# Update item price/profit data (columns H-N) across all leve sheets
# per the scheduled market-data refresh (Universalis price snapshot).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 62501092
$ws.Range("I135").Value = 25001190
$ws.Range("J135").Value = 250000600
$ws.Range("K135").Value = 225010710
$ws.Range("L135").Value = 2250005400
$ws.Range("M135").Value = -225008175
$ws.Range("N135").Value = -2250010470
$ws.Range("H138").Value = 2562.3164
$ws.Range("I138").Value = 1053.4375
$ws.Range("J138").Value = 4898.645
$ws.Range("K138").Value = 3160.3125
$ws.Range("L138").Value = 14695.935
$ws.Range("M138").Value = 1979.6875
$ws.Range("N138").Value = -24975.935

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20033.238
$ws.Range("I32").Value = 21593.164
$ws.Range("K32").Value = 21593.164
$ws.Range("M32").Value = -21306.164
$ws.Range("H61").Value = 6549.8105
$ws.Range("I61").Value = 3394.3572
$ws.Range("J61").Value = 14832.875
$ws.Range("K61").Value = 3394.3572
$ws.Range("L61").Value = 14832.875
$ws.Range("M61").Value = -3182.3572
$ws.Range("N61").Value = -15256.875
$ws.Range("H74").Value = 5587.4443
$ws.Range("I74").Value = 1953
$ws.Range("K74").Value = 1953
$ws.Range("M74").Value = -1079
$ws.Range("H77").Value = 5587.4443
$ws.Range("I77").Value = 1953
$ws.Range("K77").Value = 9765
$ws.Range("M77").Value = -5397
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H132").Value = 1712.0941
$ws.Range("I132").Value = 1328.1746
$ws.Range("J132").Value = 2811.5
$ws.Range("K132").Value = 3984.5238
$ws.Range("L132").Value = 8434.5
$ws.Range("M132").Value = -1454.5238
$ws.Range("N132").Value = -13494.5
$ws.Range("H136").Value = 6549.8105
$ws.Range("I136").Value = 3394.3572
$ws.Range("J136").Value = 14832.875
$ws.Range("K136").Value = 10183.0716
$ws.Range("L136").Value = 44498.625
$ws.Range("M136").Value = -7633.071599999999
$ws.Range("N136").Value = -49598.625

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20666.666
$ws.Range("J92").Value = 20666.666
$ws.Range("L92").Value = 20666.666
$ws.Range("N92").Value = -25658.666
$ws.Range("H103").Value = 39250
$ws.Range("J103").Value = 39250
$ws.Range("L103").Value = 39250
$ws.Range("N103").Value = -41594
$ws.Range("H118").Value = 57139.2
$ws.Range("J118").Value = 57139.2
$ws.Range("L118").Value = 57139.2
$ws.Range("N118").Value = -60453.2
$ws.Range("H134").Value = 30494.514
$ws.Range("I134").Value = 1767.1428
$ws.Range("K134").Value = 5301.428400000001
$ws.Range("M134").Value = -2766.428400000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2124.1133
$ws.Range("I31").Value = 1495.7693
$ws.Range("J31").Value = 3874.5
$ws.Range("K31").Value = 1495.7693
$ws.Range("L31").Value = 3874.5
$ws.Range("M31").Value = -1200.7693
$ws.Range("N31").Value = -4464.5
$ws.Range("H34").Value = 2124.1133
$ws.Range("I34").Value = 1495.7693
$ws.Range("J34").Value = 3874.5
$ws.Range("K34").Value = 1495.7693
$ws.Range("L34").Value = 3874.5
$ws.Range("M34").Value = -1293.7693
$ws.Range("N34").Value = -4278.5
$ws.Range("H99").Value = 3822.6365
$ws.Range("I99").Value = 3266.875
$ws.Range("J99").Value = 5304.6665
$ws.Range("K99").Value = 3266.875
$ws.Range("L99").Value = 5304.6665
$ws.Range("M99").Value = -1768.875
$ws.Range("N99").Value = -8300.666499999999
$ws.Range("H126").Value = 3822.6365
$ws.Range("I126").Value = 3266.875
$ws.Range("J126").Value = 5304.6665
$ws.Range("K126").Value = 9800.625
$ws.Range("L126").Value = 15913.9995
$ws.Range("M126").Value = -7330.625
$ws.Range("N126").Value = -20853.9995

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 50000590
$ws.Range("J5").Value = 125000480
$ws.Range("L5").Value = 375001440
$ws.Range("N5").Value = -375001664
$ws.Range("H18").Value = 343.5263
$ws.Range("I18").Value = 266.29413
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 798.88239
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -629.88239
$ws.Range("N18").Value = -3338
$ws.Range("H113").Value = 699.1395
$ws.Range("I113").Value = 697
$ws.Range("J113").Value = 708.5
$ws.Range("K113").Value = 2091
$ws.Range("L113").Value = 2125.5
$ws.Range("M113").Value = 79
$ws.Range("N113").Value = -6465.5
$ws.Range("H122").Value = 950.6316
$ws.Range("I122").Value = 495.6
$ws.Range("K122").Value = 4460.400000000001
$ws.Range("M122").Value = -2010.400000000001
$ws.Range("H135").Value = 50000590
$ws.Range("J135").Value = 125000480
$ws.Range("L135").Value = 1125004320
$ws.Range("N135").Value = -1125009390

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 26000
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5411.7
$ws.Range("I7").Value = 4901.8887
$ws.Range("K7").Value = 4901.8887
$ws.Range("M7").Value = -4789.8887
$ws.Range("H55").Value = 502.44446
$ws.Range("I55").Value = 444
$ws.Range("J55").Value = 575.5
$ws.Range("K55").Value = 444
$ws.Range("L55").Value = 575.5
$ws.Range("M55").Value = -271
$ws.Range("N55").Value = -921.5
$ws.Range("H61").Value = 640230.0600000001
$ws.Range("I61").Value = 18475.076
$ws.Range("J61").Value = 3334501.8
$ws.Range("K61").Value = 18475.076
$ws.Range("L61").Value = 3334501.8
$ws.Range("M61").Value = -18273.076
$ws.Range("N61").Value = -3334905.8
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H113").Value = 640230.0600000001
$ws.Range("I113").Value = 18475.076
$ws.Range("J113").Value = 3334501.8
$ws.Range("K113").Value = 18475.076
$ws.Range("L113").Value = 3334501.8
$ws.Range("M113").Value = -16305.076
$ws.Range("N113").Value = -3338841.8
$ws.Range("H122").Value = 6495.1562
$ws.Range("I122").Value = 6151.3955
$ws.Range("J122").Value = 7199.048
$ws.Range("K122").Value = 18454.1865
$ws.Range("L122").Value = 21597.144
$ws.Range("M122").Value = -16004.1865
$ws.Range("N122").Value = -26497.144
$ws.Range("H126").Value = 5411.7
$ws.Range("I126").Value = 4901.8887
$ws.Range("K126").Value = 14705.6661
$ws.Range("M126").Value = -12235.6661
$ws.Range("H136").Value = 3502.2058
$ws.Range("I136").Value = 1936.159
$ws.Range("J136").Value = 6373.2915
$ws.Range("K136").Value = 5808.477000000001
$ws.Range("L136").Value = 19119.8745
$ws.Range("M136").Value = -3258.477000000001
$ws.Range("N136").Value = -24219.8745

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 40062
$ws.Range("J49").Value = 40062
$ws.Range("L49").Value = 40062
$ws.Range("N49").Value = -40522
$ws.Range("H68").Value = 37635.5
$ws.Range("J68").Value = 37635.5
$ws.Range("L68").Value = 37635.5
$ws.Range("N68").Value = -39257.5
$ws.Range("H71").Value = 37635.5
$ws.Range("J71").Value = 37635.5
$ws.Range("L71").Value = 112906.5
$ws.Range("N71").Value = -121018.5
$ws.Range("H80").Value = 37801
$ws.Range("J80").Value = 37801
$ws.Range("L80").Value = 37801
$ws.Range("N80").Value = -39797
$ws.Range("H83").Value = 37801
$ws.Range("J83").Value = 37801
$ws.Range("L83").Value = 113403
$ws.Range("N83").Value = -123387

